$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 3 into row 4 (keeps styles/number formats identical, like
# copying an existing booking row to start a new one for debugging).
$ws.Range("A3:L3").Copy($ws.Range("A4:L4"))

# New booking: SNOW-003 for John Smith2.
$ws.Range("A4").Value = "SNOW-003"
$ws.Range("C4").Value = "John Smith2"

# Move the active selection like a user would after entering the new row.
$ws.Range("C7").Select() | Out-Null

# Nudge the saved window position (as captured in the workbook view).
$excel.ActiveWindow.Left = 5688
$excel.ActiveWindow.Top = 3192
